# Applies the edits described by the target diff to the active document.
$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replaceText, 2) | Out-Null
}

# Find.Execute's replacement path runs the text through a "smart quotes"
# AutoFormat simulation that turns straight apostrophes into curly ones.
# When the target text must keep a literal straight apostrophe, locate the
# range with Find (no replacement) and assign .Text directly instead - that
# path does not rewrite quote characters.
function Replace-Text-Literal($findText, $replaceText) {
    $rng = $d.Content.Duplicate
    $rng.Find.Execute($findText, $true, $false, $false, $false, $false,
                       $true, 1, $false, "", 0) | Out-Null
    if ($rng.Find.Found) {
        $rng.Text = $replaceText
    }
}

# 1. Facsimile line: drop the stray trailing digits.
Replace-Text "Facsimile: (323) 937-0958 5789465745354435" "Facsimile: (323) 937-0958"

# 2. Court / venue line.
Replace-Text "COUNTY OF  LOS ANGELES – BEVERLY HILLS COURTHOUSE " "COUNTY OF  LOS ANGELES - CENTRAL DISTRICT - UNLIMITED "

# 3. Caption party block - plaintiff name.
Replace-Text "TRAINING MATE, LLC  , an individual;" "KENNETH KLOTZLE , an individual;"

# 4. Caption party block - defendant name(s).
Replace-Text "CAROLYN KONOPA ; Does 1 TO 10, Inclusive," "EMIL JAQUIAS, JOSE ZAGO AND SUSANA ZAGO ; Does 1 TO 10, Inclusive,"

# 5. Case number.
Replace-Text "Case No. 23SMCV00406" "Case No. 22STCV06411"

# 6. Title block heading (straight apostrophe - body copy, not the footer).
Replace-Text-Literal "PLAINTIFF'S RESPONSES TO REQUESTS FOR ADMISSION, SET TWO TO PLAINTIFF CAROLYN KONOPA PROPOUNDED BY DEFENDANT TRAINING MATE, LLC  " "PLAINTIFF'S RESPONSES TO REQUEST FOR ADMISSIONS "

# 7. Propounding party.
Replace-Text "Defendant, CAROLYN KONOPA " "Defendant, EMIL JAQUIAS, JOSE ZAGO AND SUSANA ZAGO "

# 8. Responding party.
Replace-Text " Plaintiff, TRAINING MATE, LLC  " " Plaintiff, KENNETH KLOTZLE "

# 9. Preliminary "makes the following answers..." sentence.
Replace-Text "PLAINTIFF, TRAINING MATE, LLC  , (hereinafter “Responding Party”) makes the following answers and objections to the REQUESTS FOR ADMISSION, SET TWO TO PLAINTIFF CAROLYN KONOPA PROPOUNDED BY DEFENDANT TRAINING MATE, LLC  , Set One. " "PLAINTIFF, KENNETH KLOTZLE , (hereinafter “Responding Party”) makes the following answers and objections to the REQUEST FOR ADMISSIONS , Set One. "

# 10. Request No. 34 -> No. 1 heading.
Replace-Text "REQUEST FOR ADMISSION NO. 34:" "REQUEST FOR ADMISSION NO. 1:"

# 11. Body of Request No. 1 (was the Training Mate liability-release text).
Replace-Text "           On June 7, 2018, YOU checked a box on the Training Mate website indicating that YOU agreed with the Liability Release language set forth below: In exchange for participation in the activity of Group fitness organized by Training Mate LLC (`"Training Mate LLC`"), of 7825 Santa Monica Blvd, West Hollywood, California, 90046 and/or use of the property, facilities and services of Training Mate LLC, I agree for myself and (if applicable) for the members of my family, to the following: 1. I agree to observe and obey all posted rules and warnings, and further agree to follow any oral instructions or directions given by Training Mate, or the employees, representatives or agents of Training Mate. 2. I recognize that there are certain inherent risks associated with the above described activity and I assume full responsibility for personal injury to myself and (if applicable) my family members, and further release and discharge Training Mate for injury, loss or damage arising out of my or my family's use of or presence upon the facilities of Training Mate, whether caused by the fault of myself, my family, Training Mate or other third parties. 3. I agree to indemnify and defend Training Mate against all claims, causes of action, damages, judgments, costs or expenses, including attorney fees and other litigation costs, which may in any way arise from my or my family's use of or presence upon the facilities of Training Mate. 4. I agree to pay for all damages to the facilities of Training Mate caused by my or my family's negligent, reckless, or willful actions. 5. I agree Training Mate, or the employees, representatives or agents of Training Mate consent to use my photographs and audio-visual recordings to promote the services offered by Training Mate LLC in advertising publications, marketing materials, publicity, or promotion. These collective images and recordings may be used on the Training Mate website, Training Mate emails, Training Mate Social Media Outlets including Facebook, Instagram, Twitter and Pinterest. I hereby assign all rights to the Training Mate, or the employees, representatives or agents of Training Mate (collectively, the “Training Mate”), and release the Training Mate from liability that may arise from the photographs or recordings. 6. Any legal or equitable claim that may arise from participation in the above shall be resolved under California law. 6. Member is required to cancel 12 hours before class, or they will be charged a `$10 late penalty fee for that class I HAVE READ THIS DOCUMENT AND ☑ I agree with the above terms" "           Admit that you cannot establish financial responsibility as required by California Civil Code Section 3333.4, enacted by the voters of the State of California through Initiative Measure Proposition 213 and effective November 5, 1996, for the time of the accident which forms the basis of this lawsuit. 330 North Brand Blvd., Suite 900 Glendale, CA 91203-2340 Electronic Address:  LosAngelesLegal@allstate.com Telephone:  (818) 548-6381 By:  LOWELL G. HOUGHTON State Bar No. 206372 Our File No. 0609106323.1- Attorneys for Defendant(s): EMIL JAQUIAS, JOSE ZAGO and SUSANA ZAGO"

# 12. Remove Requests For Admission No. 35 through No. 51 (headings, bodies,
#     responses) - the paragraph immediately after No. 1's "RESPONSE TO
#     REQUEST FOR ADMISSION:" heading through the very last "Admit." that
#     precedes the "Dated:" signature block.
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd()
    if ($t -eq "REQUEST FOR ADMISSION NO. 1:") {
        # The "Admit." response paragraph to this request is 3 paragraphs later:
        # NO. 1: / body / RESPONSE TO REQUEST FOR ADMISSION: / Admit.
        $startPara = $i + 3
    }
    if ($t -eq "REQUEST FOR ADMISSION NO. 51:") {
        # The final "Admit." response paragraph is 3 paragraphs later.
        $endPara = $i + 3 - 1
    }
}
if ($startPara -ne $null -and $endPara -ne $null) {
    $rStart = $d.Paragraphs.Item($startPara).Range.Start
    $rEnd = $d.Paragraphs.Item($endPara).Range.End
    $d.Range($rStart, $rEnd).Delete() | Out-Null
}

# 13. Turn the remaining (now renumbered) "Admit." response into an
#     objection + admission, on two lines within the same run.
Replace-Text "           Admit. " "           Objection. This discovery request seeks attorney work product in violation of Code of Civil Procedure sections 2018.020 and 2018.030. Notwithstanding the foregoing objections and subject thereto, Responding Party responds as follows: ^lAdmit. "

# 14. Footer title (curly apostrophe - distinct run from the body heading).
#     $d.Content only covers the main body story, so reach into the
#     section's footer range explicitly.
$footer = $d.Sections.Item(1).Footers.Item(1)
$footer.Range.Find.Execute("PLAINTIFF’S RESPONSES TO REQUESTS FOR ADMISSION, SET TWO TO PLAINTIFF CAROLYN KONOPA PROPOUNDED BY DEFENDANT TRAINING MATE, LLC  ", $true, $false, $false, $false, $false,
                            $true, 1, $false, "PLAINTIFF’S RESPONSES TO REQUEST FOR ADMISSIONS ", 2) | Out-Null

Write-Output "done"
